$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "64.655.95"
$ws.Range("E2").Value = "  +0.51%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.164.20"
$ws.Range("E3").Value = "  +3.09%  "
$ws.Range("E4").Value = "  +0.23%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "563.80"
$ws.Range("E5").Value = "  +2.11%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "142.98"
$ws.Range("E6").Value = "  +2.62%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.157.90"
$ws.Range("E8").Value = "  +3.16%  "
$ws.Range("E9").Value = "  +1.64%  "
$ws.Range("E10").Value = "  +4.77%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.154"
$ws.Range("E11").Value = "  +1.07%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.467"
$ws.Range("E12").Value = "  +2.11%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "36.68"
$ws.Range("E13").Value = "  +2.96%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000222"
$ws.Range("E14").Value = "  +1.43%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.676.36"
$ws.Range("E15").Value = "  +3.37%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.738.14"
$ws.Range("E16").Value = "  +0.55%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.162.08"
$ws.Range("E17").Value = "  +3.24%  "
$ws.Range("E18").Value = "  +1.04%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "517.50"
$ws.Range("E19").Value = "  +6.16%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.87"
$ws.Range("E20").Value = "  +3.79%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.08"
$ws.Range("E21").Value = "  +3.03%  "
$ws.Range("E22").Value = "  +5.02%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.48"
$ws.Range("E23").Value = "  +4.44%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.80"
$ws.Range("E24").Value = "  +3.30%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "79.21"
$ws.Range("E25").Value = "  +1.30%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("E26").Value = "  +0.03%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.83"
$ws.Range("E27").Value = "  +13.12%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.85"
$ws.Range("E28").Value = "  +4.84%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.15"
$ws.Range("E29").Value = "  +3.35%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.00"
$ws.Range("E30").Value = "  +0.23%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "26.65"
$ws.Range("E31").Value = "  +2.89%  "
$ws.Range("E32").Value = "  -0.40%  "
$ws.Range("E33").Value = "  +2.59%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "552.40"
$ws.Range("E34").Value = "  -3.89%  "
$ws.Range("E35").Value = "  +2.90%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.38"
$ws.Range("E36").Value = "  -1.18%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "54.13"
$ws.Range("E37").Value = "  +3.45%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0434"
$ws.Range("E38").Value = "  +6.68%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0826"
$ws.Range("E39").Value = "  +4.11%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.170.17"
$ws.Range("E40").Value = "  +7.47%  "
$ws.Range("E41").Value = "  +4.11%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.79"
$ws.Range("E42").Value = "  -3.20%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.32"
$ws.Range("E43").Value = "  +1.01%  "
$ws.Range("E44").Value = "  +9.69%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.22"
$ws.Range("E45").Value = "  +6.77%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "25.21"
$ws.Range("E47").Value = "  +0.29%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "121.84"
$ws.Range("E48").Value = "  +2.69%  "
$ws.Range("E49").Value = "  +0.61%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0₃0521"
$ws.Range("E50").Value = "  -1.54%  "
$ws.Range("B51").Value = "CoreDAO"
$ws.Range("C51").Value = "https://coinranking.com/coin/HFvoXUQh4+coredao-core"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.41"
$ws.Range("E51").Value = "  +58.83%  "
